# usecase_description_강지완.xlsx
#
# docs: use case diagram, use case description 수정
#       - 자전거 대여소 리스트 조회, 대여소 상세정보 조회 use case 분리
#
# The "자전거 대여소 리스트 조회" (#4) use case used to include the
# "select a station -> show its detail screen -> delete from the detail
# screen" steps. Those are split out into a brand-new use case #5
# "자전거 대여소 상세정보 조회", and use case #4 is reworded so the delete
# flow happens directly from the list (no detail screen involved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Reword steps 3-7 of use case #4 (자전거 대여소 리스트 조회): deleting a
#    station now happens straight from the list instead of via a detail
#    screen.
# ---------------------------------------------------------------------------
$ws.Range("B31").Value = "3. 관리자가 자전거 대여소 리스트에서 삭제 버튼 선택"
$ws.Range("C32").Value = "4. 시스템은 선택된 대여소의 삭제 여부를 확인"
$ws.Range("B33").Value = "5. 관리자가 삭제를 확인"
$ws.Range("C34").Value = "6. 시스템은 해당 대여소를 DB에서 삭제하고 리스트를 갱신"
$ws.Range("B35").Value = ""
$ws.Range("C35").Value = "7. 시스템은 삭제 확인 메시지를 출력"
# C36 / C37 (8. / 9.) keep their original wording.

# The A28:A37 "4" label column is no longer a clean single merge once the
# use case was edited/split - unmerge it and drop the centering that went
# with the merged look.
$aCol = $ws.Range("A28:A37")
if ($aCol.MergeCells()) {
    $aCol.UnMerge()
}
$aCol.HorizontalAlignment = 1

# Row 37 loses the slightly-taller leftover height it used to have.
$ws.Rows.Item(37).RowHeight = 49.5

# ---------------------------------------------------------------------------
# 2) Insert the new use case #5 "자전거 대여소 상세정보 조회" right after
#    use case #4, separated by one blank row (row 38), matching the blank
#    separator rows used between every other use case block.
# ---------------------------------------------------------------------------

# Row 39: section header ("5" | "자전거 대여소 상세정보 조회")
$ws.Rows.Item(27).Copy()
$ws.Rows.Item(39).Insert()
$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "자전거 대여소 상세정보 조회"

# Row 40: "Actor action" / "System Response" column header
$ws.Rows.Item(28).Copy()
$ws.Rows.Item(40).Insert()
$ws.Range("B40").Value = "Actor action"
$ws.Range("C40").Value = "System Response"

# Row 41: "1. 시스템은 자전거 대여소 리스트 표시" (system response only)
$ws.Rows.Item(30).Copy()
$ws.Rows.Item(41).Insert()
$ws.Range("B41").Value = ""
$ws.Range("C41").Value = "1. 시스템은 자전거 대여소 리스트 표시"

# Row 42: "2. 관리자가 특정 대여소를 선택" (actor action only)
$ws.Rows.Item(29).Copy()
$ws.Rows.Item(42).Insert()
$ws.Range("B42").Value = "2. 관리자가 특정 대여소를 선택"
$ws.Range("C42").Value = ""

# Row 43: "3. 시스템은 선택된 대여소의 상세정보 화면을 표시" (system response only)
$ws.Rows.Item(30).Copy()
$ws.Rows.Item(43).Insert()
$ws.Range("B43").Value = ""
$ws.Range("C43").Value = "3. 시스템은 선택된 대여소의 상세정보 화면을 표시"

# Row heights for the new block (header rows match the other section
# headers; step rows match the taller wrapped-text rows used for step
# content elsewhere in the sheet).
$ws.Rows.Item(39).RowHeight = 24.95
$ws.Rows.Item(40).RowHeight = 24.95
$ws.Rows.Item(41).RowHeight = 50.1
$ws.Rows.Item(42).RowHeight = 50.1
$ws.Rows.Item(43).RowHeight = 50.1

# Re-create the merges for the new block: the number/title row and the
# "A" label column spanning all four Actor/System rows.
$ws.Range("B39:C39").Merge()
$ws.Range("A40:A43").Merge()
$ws.Range("A40:A43").HorizontalAlignment = -4108

Write-Host "Done"
